$d = $word.ActiveDocument

# --- Fix 1: "NUMERO 23" paragraph currently has a proofing-error split
# (separate runs around a gramStart/gramEnd proofErr pair). Rebuild it as
# a single clean paragraph containing one run "NUMERO 23".
$pSchindler = $d.Paragraphs.Item(3)   # "A LISTA DE SCHINDLER"
$pSchindler.Range.InsertParagraphAfter()
$pNumeroNew = $d.Paragraphs.Item(4)
$pNumeroNew.Range.Text = "NUMERO 23"

# The old "NUMERO 23" paragraph (with the proofErr markup) got pushed down
# to index 5; remove it entirely.
$pNumeroOld = $d.Paragraphs.Item(5)
$pNumeroOld.Range.Delete()

# --- Fix 2: add a new movie "VINGADORES" as its own paragraph right after
# "DEBI E LOIDE".
$pDebi = $d.Paragraphs.Item(6)        # "DEBI E LOIDE"
$pDebi.Range.InsertParagraphAfter()
$pVingadores = $d.Paragraphs.Item(7)
$pVingadores.Range.Text = "VINGADORES"
